# Rename the three inline picture shapes (wp:docPr/@name and pic:cNvPr/@name)
# found in the document's header/footer parts:
#   footer1.xml : id="3" image "image1.png" -> "image2.png"
#   footer2.xml : id="2" image "image1.png" -> "image2.png"
#   header1.xml : id="1" image "image2.jpg" -> "image1.jpg"
#
# InlineShape objects expose no writable "Name" property in the Word object
# model (this matches real Word/VBA - InlineShape has no .Name), so the
# rename is applied by round-tripping the package through
# Document.WordOpenXML and patching the three distinctive XML fragments
# that contain both the wp:docPr and pic:cNvPr elements for each picture.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$old1 = 'id="3" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
$new1 = 'id="3" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'

$old2 = 'id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
$new2 = 'id="2" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'

$old3 = 'id="1" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'
$new3 = 'id="1" name="image1.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>'

$xml = $xml.Replace($old1, $new1)
$xml = $xml.Replace($old2, $new2)
$xml = $xml.Replace($old3, $new3)

$d.WordOpenXML = $xml
